# Applies the "UML Diagram" heading edits described by the commit:
#   - Capitalises "class"/"interface" to "Class"/"Interface" in several
#     UML diagram headings (CardsDriver, LivesIn, Player, MainMenuGUI).
#   - Normalises the LivesGUI / Card headings into single runs.
#   - Renames the final "UML Diagram" heading to "VOPC Diagram" and
#     relocates the "_GoBack" bookmark into that heading.

$d = $word.ActiveDocument

$ENDASH = [char]8211

# Helper: wrap a fragment of run-level OOXML (<w:r>...</w:r>, optionally
# interspersed with <w:proofErr/>, <w:bookmarkStart/>, etc.) in the
# pkg:package envelope that Range.InsertXML expects, and replace the
# supplied Range's contents with it.
function Set-ParagraphRuns($range, [string]$innerXml) {
    $full = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($full)
}

# ---------------------------------------------------------------------
# 1. "UML Diagram - CardsDriver class" -> "... CardsDriver Class"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(73)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">UML </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">Diagram - CardsDriver </w:t></w:r>' +
         '<w:r><w:t>C</w:t></w:r>' +
         '<w:r><w:t>lass</w:t></w:r>'
Set-ParagraphRuns $r $inner

# ---------------------------------------------------------------------
# 2. "UML Diagram - LivesIn interface" -> "... LivesIn Interface"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(101)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r><w:t xml:space="preserve">UML Diagram - </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>LivesIn</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:r><w:t>I</w:t></w:r>' +
         '<w:r><w:t>nterface</w:t></w:r>'
Set-ParagraphRuns $r $inner

# ---------------------------------------------------------------------
# 3. "UML Diagram - Player class" -> "... Player Class"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(106)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">UML Diagram - Player </w:t></w:r>' +
         '<w:r><w:t>C</w:t></w:r>' +
         '<w:r><w:t>lass</w:t></w:r>'
Set-ParagraphRuns $r $inner

# ---------------------------------------------------------------------
# 4. "UML Diagram - MainMenuGUI class" -> "... MainMenuGUI Class"
#    (the dash here is an en-dash)
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(111)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">UML Diagram ' + $ENDASH + ' MainMenuGUI </w:t></w:r>' +
         '<w:r><w:t>C</w:t></w:r>' +
         '<w:r><w:t>lass</w:t></w:r>'
Set-ParagraphRuns $r $inner

# ---------------------------------------------------------------------
# 5. "UML Diagram - LivesGUI Class" - collapse into a single run and
#    drop the stray "_GoBack" bookmark that used to sit here (it moves
#    to the final "VOPC Diagram" heading below).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$p = $d.Paragraphs.Item(132)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r><w:lastRenderedPageBreak/><w:t>UML Diagram ' + $ENDASH + ' LivesGUI Class</w:t></w:r>'
Set-ParagraphRuns $r $inner

# ---------------------------------------------------------------------
# 6. "UML Diagram - Card Class" - collapse into a single run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(134)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r><w:lastRenderedPageBreak/><w:t>UML Diagram - Card Class</w:t></w:r>'
Set-ParagraphRuns $r $inner

# ---------------------------------------------------------------------
# 7. "UML Diagram" -> "VOPC Diagram", with the "_GoBack" bookmark
#    re-inserted right after "VOPC".
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(143)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r><w:lastRenderedPageBreak/><w:t>VOPC</w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
         '<w:bookmarkEnd w:id="0"/>' +
         '<w:r><w:t xml:space="preserve"> Diagram</w:t></w:r>'
Set-ParagraphRuns $r $inner
